$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 4006.389
$ws.Range("I40").Value = 2724.25
$ws.Range("J40").Value = 4372.7144
$ws.Range("K40").Value = 2724.25
$ws.Range("L40").Value = 4372.7144
$ws.Range("M40").Value = -2549.25
$ws.Range("N40").Value = -4722.7144
# Row 76
$ws.Range("H76").Value = 4998.5
$ws.Range("J76").Value = 4998
$ws.Range("L76").Value = 4998
$ws.Range("N76").Value = -5628
# Row 79
$ws.Range("H79").Value = 4998.5
$ws.Range("J79").Value = 4998
$ws.Range("L79").Value = 4998
$ws.Range("N79").Value = -7182
# Row 80
$ws.Range("H80").Value = 32343.525
$ws.Range("J80").Value = 872.26666
$ws.Range("L80").Value = 2616.79998
$ws.Range("N80").Value = -4612.79998
# Row 83
$ws.Range("H83").Value = 32343.525
$ws.Range("J83").Value = 872.26666
$ws.Range("L83").Value = 7850.39994
$ws.Range("N83").Value = -17834.39994
# Row 86
$ws.Range("H86").Value = 69337.8
$ws.Range("I86").Value = 103087.1
$ws.Range("J86").Value = 1839.2
$ws.Range("K86").Value = 103087.1
$ws.Range("L86").Value = 1839.2
$ws.Range("M86").Value = -101964.1
$ws.Range("N86").Value = -4085.2
# Row 89
$ws.Range("H89").Value = 69337.8
$ws.Range("I89").Value = 103087.1
$ws.Range("J89").Value = 1839.2
$ws.Range("K89").Value = 515435.5
$ws.Range("L89").Value = 9196
$ws.Range("M89").Value = -509819.5
$ws.Range("N89").Value = -20428
# Row 111
$ws.Range("H111").Value = 943.2
$ws.Range("I111").Value = 1012.44446
$ws.Range("K111").Value = 3037.33338
$ws.Range("M111").Value = 29.66661999999997
# Row 115
$ws.Range("H115").Value = 1083.6666
$ws.Range("I115").Value = 261
$ws.Range("K115").Value = 783
$ws.Range("M115").Value = 784
# Row 116
$ws.Range("H116").Value = 5265.4146
$ws.Range("I116").Value = 5117.8184
$ws.Range("K116").Value = 5117.8184
$ws.Range("M116").Value = -1675.8184
# Row 125
$ws.Range("H125").Value = 1374.75
$ws.Range("I125").Value = 1124
$ws.Range("K125").Value = 10116
$ws.Range("M125").Value = -7656
# Row 137
$ws.Range("H137").Value = 5479.2666
$ws.Range("I137").Value = 1288.4286
$ws.Range("J137").Value = 9146.25
$ws.Range("K137").Value = 3865.2858
$ws.Range("L137").Value = 27438.75
$ws.Range("M137").Value = -1315.2858
$ws.Range("N137").Value = -32538.75
# Row 141
$ws.Range("H141").Value = 5833.222
$ws.Range("I141").Value = 6287.375
$ws.Range("K141").Value = 18862.125
$ws.Range("M141").Value = -13682.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 15000
$ws.Range("J24").Value = 15000
$ws.Range("L24").Value = 15000
$ws.Range("N24").Value = -15748
# Row 32
$ws.Range("H32").Value = 70405.88
$ws.Range("I32").Value = 86992.30499999999
$ws.Range("K32").Value = 86992.30499999999
$ws.Range("M32").Value = -86705.30499999999
# Row 61
$ws.Range("H61").Value = 4716.6665
$ws.Range("I61").Value = 4700
$ws.Range("J61").Value = 4800
$ws.Range("K61").Value = 4700
$ws.Range("L61").Value = 4800
$ws.Range("M61").Value = -4488
$ws.Range("N61").Value = -5224
# Row 88
$ws.Range("H88").Value = 1441.3043
$ws.Range("I88").Value = 991
$ws.Range("J88").Value = 1681.4667
$ws.Range("K88").Value = 991
$ws.Range("L88").Value = 1681.4667
$ws.Range("M88").Value = -585
$ws.Range("N88").Value = -2493.4667
# Row 91
$ws.Range("H91").Value = 1441.3043
$ws.Range("I91").Value = 991
$ws.Range("J91").Value = 1681.4667
$ws.Range("K91").Value = 991
$ws.Range("L91").Value = 1681.4667
$ws.Range("M91").Value = 413
$ws.Range("N91").Value = -4489.4667
# Row 97
$ws.Range("H97").Value = 5235.0938
$ws.Range("I97").Value = 6437.35
$ws.Range("J97").Value = 3231.3333
$ws.Range("K97").Value = 6437.35
$ws.Range("L97").Value = 3231.3333
$ws.Range("M97").Value = -5941.35
$ws.Range("N97").Value = -4223.3333
# Row 100
$ws.Range("H100").Value = 15000
$ws.Range("J100").Value = 15000
$ws.Range("L100").Value = 15000
$ws.Range("N100").Value = -17164
# Row 102
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").Value = $null
# Row 122
$ws.Range("H122").Value = 2250.6
$ws.Range("I122").Value = 2168.5
$ws.Range("J122").Value = 2373.75
$ws.Range("K122").Value = 6505.5
$ws.Range("L122").Value = 7121.25
$ws.Range("M122").Value = -4055.5
$ws.Range("N122").Value = -12021.25
# Row 132
$ws.Range("H132").Value = 504443.5
$ws.Range("I132").Value = 999999
$ws.Range("J132").Value = 8888
$ws.Range("K132").Value = 2999997
$ws.Range("L132").Value = 26664
$ws.Range("M132").Value = -2997467
$ws.Range("N132").Value = -31724
# Row 136
$ws.Range("H136").Value = 4716.6665
$ws.Range("I136").Value = 4700
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 14100
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -11550
$ws.Range("N136").Value = -19500

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("N20").Value = $null

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 15
$ws.Range("H15").Value = 400
$ws.Range("I15").Value = 400
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 400
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -230
$ws.Range("N15").Value = $null
# Row 22
$ws.Range("H22").Value = 570.1429000000001
$ws.Range("I22").Value = 341.5
$ws.Range("K22").Value = 341.5
$ws.Range("M22").Value = 8.5
# Row 31
$ws.Range("H31").Value = 4336.136
$ws.Range("I31").Value = 3036.375
$ws.Range("J31").Value = 7802.1665
$ws.Range("K31").Value = 3036.375
$ws.Range("L31").Value = 7802.1665
$ws.Range("M31").Value = -2741.375
$ws.Range("N31").Value = -8392.166499999999
# Row 34
$ws.Range("H34").Value = 4336.136
$ws.Range("I34").Value = 3036.375
$ws.Range("J34").Value = 7802.1665
$ws.Range("K34").Value = 3036.375
$ws.Range("L34").Value = 7802.1665
$ws.Range("M34").Value = -2834.375
$ws.Range("N34").Value = -8206.166499999999
# Row 94
$ws.Range("H94").Value = 3247.3333
$ws.Range("J94").Value = 3182.4285
$ws.Range("L94").Value = 3182.4285
$ws.Range("N94").Value = -4084.4285

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 35.76923
$ws.Range("I2").Value = 35.3
$ws.Range("J2").Value = 37.333332
$ws.Range("K2").Value = 211.8
$ws.Range("L2").Value = 223.999992
$ws.Range("M2").Value = -98.79999999999998
$ws.Range("N2").Value = -449.999992
# Row 107
$ws.Range("H107").Value = 873.2
$ws.Range("I107").Value = 1885
$ws.Range("K107").Value = 5655
$ws.Range("M107").Value = -3735
# Row 131
$ws.Range("H131").Value = 1965622.5
$ws.Range("J131").Value = 2783931
$ws.Range("L131").Value = 8351793
$ws.Range("N131").Value = -8361873
# Row 132
$ws.Range("H132").Value = 1495.5
$ws.Range("I132").Value = 1632.6666
$ws.Range("J132").Value = 1358.3334
$ws.Range("K132").Value = 14693.9994
$ws.Range("L132").Value = 12225.0006
$ws.Range("M132").Value = -12163.9994
$ws.Range("N132").Value = -17285.0006

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 3247
$ws.Range("I80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("M80").Value = $null
# Row 83
$ws.Range("H83").Value = 3247
$ws.Range("I83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("M83").Value = $null
# Row 126
$ws.Range("H126").Value = 10994.5
$ws.Range("I126").Value = 6326
$ws.Range("K126").Value = 18978
$ws.Range("M126").Value = -16508
# Row 134
$ws.Range("H134").Value = 40000
$ws.Range("J134").Value = 40000
$ws.Range("L134").Value = 120000
$ws.Range("N134").Value = -125070

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3690.2727
$ws.Range("I40").Value = 3210.3333
$ws.Range("J40").Value = 5850
$ws.Range("K40").Value = 3210.3333
$ws.Range("L40").Value = 5850
$ws.Range("M40").Value = -3074.3333
$ws.Range("N40").Value = -6122
# Row 68
$ws.Range("H68").Value = 3202.923
$ws.Range("I68").Value = 2862.3333
$ws.Range("K68").Value = 2862.3333
$ws.Range("M68").Value = -2113.3333
# Row 71
$ws.Range("H71").Value = 3202.923
$ws.Range("I71").Value = 2862.3333
$ws.Range("K71").Value = 14311.6665
$ws.Range("M71").Value = -10567.6665
# Row 125
$ws.Range("H125").Value = 74749.5
$ws.Range("J125").Value = 74749.5
$ws.Range("L125").Value = 74749.5
$ws.Range("N125").Value = -84589.5
# Row 132
$ws.Range("H132").Value = 39621.375
$ws.Range("I132").Value = 39621.375
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 118864.125
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -116334.125
$ws.Range("N132").Value = $null
# Row 135
$ws.Range("H135").Value = 85899
$ws.Range("J135").Value = 85899
$ws.Range("L135").Value = 85899
$ws.Range("N135").Value = -96039

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 264744.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 264744.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 264744.75
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -265992.75
# Row 65
$ws.Range("H65").Value = 264744.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 264744.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 1323723.75
$ws.Range("M65").Value = $null
$ws.Range("N65").Value = -1329963.75
# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
# Row 121
$ws.Range("H121").Value = 44999
$ws.Range("J121").Value = 44999
$ws.Range("L121").Value = 44999
$ws.Range("N121").Value = -48493


Write-Output "Applied all profit sheet updates."
